$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

# The remark / code-table cells used "NN label" (space separated); the
# author reformatted both lists to use "NN:label" (colon separated).
$ws.Range("G12").Value = "共用代碼檔(CustRelationType)`n01:本人`n02:配偶`n03:祖(外祖)父母`n04:父母`n05:兄弟姊妹`n06:子女`n07:孫(外孫)子女`n08:有控制與從屬關係`n09:相互投資關係`n10:董事長`n11:董事`n12:監察人`n99:其他"

$ws.Range("G13").Value = "1:持股比例`n2:被持股比例`n3:持有股份`n4:出資額`n5:關係人`n9:其它"

# Reflect the scrolled/selected state the author ended up in after editing.
$excel.ActiveWindow.ScrollRow = 12
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G14").Select()
